# Spring_Chinook_Habitat_Quality_PROTECTION_Methow_Entiat_Wenatchee_OKanogan.xlsx
# "Minor updates to data being read in, etc."
#
# 1) Two reaches were dropped from the source data pull and their rows
#    removed (remaining rows shift up):
#      - "Entiat River Potato 07"
#      - "Methow River Thompson 09"
# 2) A handful of remaining reaches got refreshed Riparian-CanopyCover /
#    Riparian-Disturbance scores (P, Q), which cascade into the derived
#    Riparian_Mean (R), HQ_Sum (T), HQ_Pct (U), HQ_Score_Restoration (V)
#    and HQ_Score_Protection (W) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-RowByReachName($name) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cellValue = $ws.Cells.Item($r, 1).Value()
        if ($cellValue -eq $name) {
            return $r
        }
    }
    return -1
}

# --- 1) Remove the two dropped reaches -------------------------------------
$rowsToDelete = @("Entiat River Potato 07", "Methow River Thompson 09")

foreach ($reachName in $rowsToDelete) {
    $rowIndex = Find-RowByReachName $reachName
    if ($rowIndex -gt 0) {
        $ws.Rows.Item($rowIndex).Delete()
    }
}

# --- 2) Refresh scores for the reaches whose inputs changed -----------------
# Each entry: ReachName -> updated column values (only the columns that changed)
$updates = @{
    "Entiat River Lake 03" = @{ P = 3; Q = 3; R = 3; T = 39; U = 0.8666666666666667; V = 3; W = 3 }
    "Entiat River Lake 04" = @{ P = 3;        R = 4; T = 32; U = 0.7111111111111111 }
    "Entiat River Lake 05" = @{ P = 3;        R = 4; T = 40; U = 0.8888888888888888; V = 3; W = 3 }
    "Entiat River Lake 06" = @{ P = 3;        R = 3; T = 37; U = 0.8222222222222222 }
    "Methow River Rattlesnake 05" = @{ P = 3; R = 3; T = 37; U = 0.8222222222222222 }
    "Methow River Thompson 07" = @{ P = 1; Q = 1; R = 1; T = 32; U = 0.7111111111111111 }
    "Methow River Thompson 08" = @{ P = 1; Q = 1; R = 1; T = 32; U = 0.7111111111111111 }
    "Nason Creek Lower 01" = @{ P = 3; Q = 3; R = 3; T = 36; U = 0.8 }
    "Nason Creek Lower 03" = @{ P = 3; Q = 3; R = 3; T = 34; U = 0.7555555555555555; V = 5 }
}

$colNumbers = @{ P = 16; Q = 17; R = 18; S = 19; T = 20; U = 21; V = 22; W = 23 }

foreach ($reachName in $updates.Keys) {
    $rowIndex = Find-RowByReachName $reachName
    if ($rowIndex -gt 0) {
        $colUpdates = $updates[$reachName]
        foreach ($colLetter in $colUpdates.Keys) {
            $colIndex = $colNumbers[$colLetter]
            $ws.Cells.Item($rowIndex, $colIndex).Value = $colUpdates[$colLetter]
        }
    }
}
